$d = $word.ActiveDocument

# 1. Append the page marker "[༦༤ན]" to the very end of the final paragraph's
#    text (right before the now-to-be-removed footnote 130 reference), and
#    2. delete the empty footnote 130 (this also removes its reference run
#    from the body).
$d.Content.Find.Execute("རྫོགས་སོ། །", $false, $false, $false, $false, $false, `
    $true, 1, $false, "རྫོགས་སོ། །[༦༤ན]", 2) | Out-Null

$footnotes = $d.Footnotes
$fn130 = $footnotes.Item($footnotes.Count)
$fn130.Delete()

# 3. Footnote 25 ("...པེ་ཅིན།aa") -> drop the stray trailing "aa".
$fn25 = $d.Footnotes.Item(5)
$r25 = $fn25.Range
$r25.MoveStart(1, 31)
$r25.Delete()

# 4. Footnote 30 (just "།") -> replace with the real note text.
$fn30 = $d.Footnotes.Item(10)
$fn30.Range.Text = "ཟས་དང་ནོར། ཞེས་པར་མ་གཞན་ནང་མེད།"
